{"js": "// Apply the textual corrections described by the commit diff.\n// (Style-id renames such as \"PargrafodaLista\" -> \"ListParagraph\" and the\n// namespace/latentStyles cleanup in the diff are artifacts of Word\n// re-serialising the package on save and are not reachable through the\n// Word JS / COM object models, so only the genuine content edits below\n// are reproduced.)\n\nconst body = context.document.body;\n\n// 1) \"... cet ensemble de cartes dans des diff\u00e9rentes moments.\"\n//    -> \"... cet ensemble de cartes dans des diff\u00e9rents moments.\"\nlet results = body.search(\"dans des diff\u00e9rentes moments.\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"dans des diff\u00e9rents moments.\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) \" et des autres attributs secondaires qu\u2019assurent le d\u00e9ro\"\n//    -> \" et des autres attributs secondaires qui assurent le d\u00e9ro\"\nresults = body.search(\"autres attributs secondaires qu\\u2019assurent le d\\u00e9ro\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"autres attributs secondaires qui assurent le d\\u00e9ro\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3) \"en focntion de l\u2019intelligence\" -> \"en fonction de l\u2019intelligence\"\nresults = body.search(\"en focntion de l\\u2019intelligence\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"en fonction de l\\u2019intelligence\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 4) \", s\u2019il est plus petit que le limite, il jouera selon \"\n//    -> \", s\u2019il est plus petit que la limite, il jouera selon \"\n//    plus a collapsed \"_GoBack\" bookmark right after \"la\" (marks the\n//    last edit position, exactly like Word does automatically).\nresults = body.search(\"que le limite\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"que la limite\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\nresults = body.search(\" limite, il jouera selon\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  const splitPoint = results.items[0].getRange(Word.RangeLocation.start);\n  await context.sync();\n  splitPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// 5) \"Note auto-attribu\u00e9e : <tab>17 /20\" -> \"... 17,5 /20\"\nresults = body.search(\"17\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < results.items.length; i++) {\n  const match = results.items[i];\n  const para = match.paragraphs.getFirst();\n  para.load(\"text\");\n  await context.sync();\n  if (para.text.indexOf(\"Note auto-attribu\\u00e9e\") === 0) {\n    match.insertText(\",5\", Word.InsertLocation.after);\n    await context.sync();\n  }\n}\n", "ps1": "# Apply the textual corrections described by the commit diff.\n# (Style-id renames such as \"PargrafodaLista\" -> \"ListParagraph\" and the\n# namespace/latentStyles cleanup seen in the diff are artifacts of Word\n# re-serialising the package on save - they are not reachable through the\n# Word COM object model, so only the genuine content edits below are\n# reproduced.)\n\n$d = $word.ActiveDocument\n\n# 1) \"... cet ensemble de cartes dans des diff\u00e9rentes moments.\"\n#    -> \"... cet ensemble de cartes dans des diff\u00e9rents moments.\"\n$r1 = $d.Content\n$r1.Find.Execute(\"dans des diff\u00e9rentes moments.\", $false, $false, $false, $false, $false, $true, 1, $false, \"dans des diff\u00e9rents moments.\", 1) | Out-Null\n\n# 2) \" et des autres attributs secondaires qu\u2019assurent le d\u00e9ro\"\n#    -> \" et des autres attributs secondaires qui assurent le d\u00e9ro\"\n$r2 = $d.Content\n$r2.Find.Execute(\"autres attributs secondaires qu\u2019assurent le d\u00e9ro\", $false, $false, $false, $false, $false, $true, 1, $false, \"autres attributs secondaires qui assurent le d\u00e9ro\", 1) | Out-Null\n\n# 3) \"en focntion de l\u2019intelligence\" -> \"en fonction de l\u2019intelligence\"\n$r3 = $d.Content\n$r3.Find.Execute(\"en focntion de l\u2019intelligence\", $false, $false, $false, $false, $false, $true, 1, $false, \"en fonction de l\u2019intelligence\", 1) | Out-Null\n\n# 4) \", s\u2019il est plus petit que le limite, il jouera selon \"\n#    -> \", s\u2019il est plus petit que la limite, il jouera selon \"\n#    plus a collapsed \"_GoBack\" bookmark right after \"la\" (marks the\n#    last edit position, exactly like Word does automatically).\n$r4 = $d.Content\n$r4.Find.Execute(\"que le limite\", $false, $false, $false, $false, $false, $true, 1, $false, \"que la limite\", 1) | Out-Null\n\n$r4b = $d.Content\n$r4b.Find.Execute(\" limite, il jouera selon\") | Out-Null\n$splitPoint = $r4b.Duplicate\n$splitPoint.Collapse(1)\n$d.Bookmarks.Add(\"_GoBack\", $splitPoint) | Out-Null\n\n# 5) \"Note auto-attribu\u00e9e : <tab>17 /20\" -> \"... 17,5 /20\"\n$r5 = $d.Content\n$r5.Find.Execute(\"17 /20\") | Out-Null\n$r5.Text = \"17,5 /20\"\n"}
